# Move the second ("O'Hare All Airlines / Denver All Airlines") table, which
# currently sits stacked below the first table (rows 13-23 on the "Graph"
# sheet), so that it instead sits beside the first table (columns F:I,
# rows 1-11), then remove the now-empty leftover rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Graph")

# Cut the second table's data columns (B:E, rows 13-23 - header + 10 years of
# data) and paste it starting at F1, right next to the first table.
$ws.Range("B13:E23").Cut($ws.Range("F1"))

# The old table's A column (years + header) and the rest of the rows are no
# longer needed now that the data lives next to the first table - remove the
# whole rows.
$ws.Range("A13:A23").EntireRow.Delete()

# Re-select the area where the old rows used to be (matches the selection
# left behind in Excel after the row delete).
$ws.Range("A13:A24").Select()

# The columns were auto-fit by Excel after the move; reproduce the resulting
# widths as closely as the host allows.
$ws.Range("B1").ColumnWidth = 27.736979166666668
$ws.Range("C1").ColumnWidth = 29.877604166666668
$ws.Range("F1").ColumnWidth = 31.736979166666668
$ws.Range("G1").ColumnWidth = 33.877604166666664
$ws.Range("H1").ColumnWidth = 32.166666666666664
$ws.Range("I1").ColumnWidth = 34.307291666666664
